$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 273
$ws.Range("I9").Value = 150.4
$ws.Range("J9").Value = 375.16666
$ws.Range("K9").Value = 150.4
$ws.Range("L9").Value = 375.16666
$ws.Range("M9").Value = 18.59999999999999
$ws.Range("N9").Value = -713.16666

$ws.Range("H55").Value = 406.2
$ws.Range("I55").Value = 200.14285
$ws.Range("J55").Value = 887
$ws.Range("K55").Value = 200.14285
$ws.Range("L55").Value = 887
$ws.Range("M55").Value = 13.85714999999999
$ws.Range("N55").Value = -1315

$ws.Range("H74").Value = 7378.3
$ws.Range("I74").Value = 7387.1113
$ws.Range("K74").Value = 7387.1113
$ws.Range("M74").Value = -6451.1113

$ws.Range("H77").Value = 7378.3
$ws.Range("I77").Value = 7387.1113
$ws.Range("K77").Value = 36935.5565
$ws.Range("M77").Value = -32255.5565

$ws.Range("H80").Value = 995
$ws.Range("I80").Value = 1147.1052
$ws.Range("K80").Value = 3441.3156
$ws.Range("M80").Value = -2443.3156

$ws.Range("H83").Value = 995
$ws.Range("I83").Value = 1147.1052
$ws.Range("K83").Value = 10323.9468
$ws.Range("M83").Value = -5331.9468

$ws.Range("H86").Value = 3431.3333
$ws.Range("I86").Value = 3476.6
$ws.Range("J86").Value = 3408.7
$ws.Range("K86").Value = 3476.6
$ws.Range("L86").Value = 3408.7
$ws.Range("M86").Value = -2353.6
$ws.Range("N86").Value = -5654.7

$ws.Range("H88").Value = 6589.926
$ws.Range("I88").Value = 3446.8333
$ws.Range("J88").Value = 9104.4
$ws.Range("K88").Value = 3446.8333
$ws.Range("L88").Value = 9104.4
$ws.Range("M88").Value = -3040.8333
$ws.Range("N88").Value = -9916.4

$ws.Range("H89").Value = 3431.3333
$ws.Range("I89").Value = 3476.6
$ws.Range("J89").Value = 3408.7
$ws.Range("K89").Value = 17383
$ws.Range("L89").Value = 17043.5
$ws.Range("M89").Value = -11767
$ws.Range("N89").Value = -28275.5

$ws.Range("H91").Value = 6589.926
$ws.Range("I91").Value = 3446.8333
$ws.Range("J91").Value = 9104.4
$ws.Range("K91").Value = 3446.8333
$ws.Range("L91").Value = 9104.4
$ws.Range("M91").Value = -2042.8333
$ws.Range("N91").Value = -11912.4

$ws.Range("H125").Value = 2756.125
$ws.Range("I125").Value = 2999
$ws.Range("J125").Value = 2721.4285
$ws.Range("K125").Value = 26991
$ws.Range("L125").Value = 24492.8565
$ws.Range("M125").Value = -24531
$ws.Range("N125").Value = -29412.8565

$ws.Range("H137").Value = 5287.5947
$ws.Range("I137").Value = 5294.5835
$ws.Range("K137").Value = 15883.7505
$ws.Range("M137").Value = -13333.7505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2369.448
$ws.Range("I32").Value = 1449.6517
$ws.Range("K32").Value = 1449.6517
$ws.Range("M32").Value = -1162.6517

$ws.Range("H61").Value = 2134.75
$ws.Range("I61").Value = 2065.182
$ws.Range("K61").Value = 2065.182
$ws.Range("M61").Value = -1853.182

$ws.Range("H63").Value = 2204.88
$ws.Range("I63").Value = 1588.8572
$ws.Range("J63").Value = 2988.9092
$ws.Range("K63").Value = 1588.8572
$ws.Range("L63").Value = 2988.9092
$ws.Range("M63").Value = -902.8571999999999
$ws.Range("N63").Value = -4360.9092

$ws.Range("H66").Value = 2204.88
$ws.Range("I66").Value = 1588.8572
$ws.Range("J66").Value = 2988.9092
$ws.Range("K66").Value = 7944.286
$ws.Range("L66").Value = 14944.546
$ws.Range("M66").Value = -4512.286
$ws.Range("N66").Value = -21808.546

$ws.Range("H88").Value = 2183.1
$ws.Range("J88").Value = 2316.889
$ws.Range("L88").Value = 2316.889
$ws.Range("N88").Value = -3128.889

$ws.Range("H91").Value = 2183.1
$ws.Range("J91").Value = 2316.889
$ws.Range("L91").Value = 2316.889
$ws.Range("N91").Value = -5124.889

$ws.Range("H122").Value = 7857.082
$ws.Range("I122").Value = 6347.375
$ws.Range("K122").Value = 19042.125
$ws.Range("M122").Value = -16592.125

$ws.Range("H136").Value = 2134.75
$ws.Range("I136").Value = 2065.182
$ws.Range("K136").Value = 6195.545999999999
$ws.Range("M136").Value = -3645.545999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2767.04
$ws.Range("I31").Value = 2248.9048
$ws.Range("J31").Value = 5487.25
$ws.Range("K31").Value = 2248.9048
$ws.Range("L31").Value = 5487.25
$ws.Range("M31").Value = -1953.9048
$ws.Range("N31").Value = -6077.25

$ws.Range("H34").Value = 2767.04
$ws.Range("I34").Value = 2248.9048
$ws.Range("J34").Value = 5487.25
$ws.Range("K34").Value = 2248.9048
$ws.Range("L34").Value = 5487.25
$ws.Range("M34").Value = -2046.9048
$ws.Range("N34").Value = -5891.25

$ws.Range("H134").Value = 2983.3333
$ws.Range("I134").Value = 2471.6487
$ws.Range("K134").Value = 7414.946100000001
$ws.Range("M134").Value = -4879.946100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4897.3335
$ws.Range("J39").Value = 5446
$ws.Range("L39").Value = 16338
$ws.Range("N39").Value = -16926

$ws.Range("H62").Value = 9224.817999999999
$ws.Range("I62").Value = 7353.2856
$ws.Range("K62").Value = 22059.8568
$ws.Range("M62").Value = -21373.8568

$ws.Range("H65").Value = 9224.817999999999
$ws.Range("I65").Value = 7353.2856
$ws.Range("K65").Value = 66179.5704
$ws.Range("M65").Value = -62747.5704

$ws.Range("H68").Value = 1844.1538
$ws.Range("J68").Value = 2097.4
$ws.Range("L68").Value = 6292.200000000001
$ws.Range("N68").Value = -7914.200000000001

$ws.Range("H71").Value = 1844.1538
$ws.Range("J71").Value = 2097.4
$ws.Range("L71").Value = 18876.6
$ws.Range("N71").Value = -26988.6

$ws.Range("H129").Value = 13368650
$ws.Range("I129").Value = 11905704
$ws.Range("J129").Value = 15876558
$ws.Range("K129").Value = 35717112
$ws.Range("L129").Value = 47629674
$ws.Range("M129").Value = -35712112
$ws.Range("N129").Value = -47639674

$ws.Range("H131").Value = 3924398.2
$ws.Range("I131").Value = 10939237
$ws.Range("K131").Value = 32817711
$ws.Range("M131").Value = -32812671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6647.75
$ws.Range("J80").Value = 3646.5
$ws.Range("L80").Value = 3646.5
$ws.Range("N80").Value = -5642.5

$ws.Range("H83").Value = 6647.75
$ws.Range("J83").Value = 3646.5
$ws.Range("L83").Value = 18232.5
$ws.Range("N83").Value = -28216.5

$ws.Range("H122").Value = 10501
$ws.Range("I122").Value = 10501
$ws.Range("K122").Value = 31503
$ws.Range("M122").Value = -29053

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 370.33334
$ws.Range("I55").Value = 133.85715
$ws.Range("J55").Value = 577.25
$ws.Range("K55").Value = 133.85715
$ws.Range("L55").Value = 577.25
$ws.Range("M55").Value = 39.14285000000001
$ws.Range("N55").Value = -923.25

$ws.Range("H82").Value = 862.46375
$ws.Range("I82").Value = 787.7451
$ws.Range("J82").Value = 1074.1666
$ws.Range("K82").Value = 787.7451
$ws.Range("L82").Value = 1074.1666
$ws.Range("M82").Value = -426.7451
$ws.Range("N82").Value = -1796.1666

$ws.Range("H85").Value = 862.46375
$ws.Range("I85").Value = 787.7451
$ws.Range("J85").Value = 1074.1666
$ws.Range("K85").Value = 787.7451
$ws.Range("L85").Value = 1074.1666
$ws.Range("M85").Value = 460.2549
$ws.Range("N85").Value = -3570.1666

$ws.Range("H132").Value = 2916.2144
$ws.Range("I132").Value = 2652.8
$ws.Range("K132").Value = 7958.400000000001
$ws.Range("M132").Value = -5428.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 25003.5
$ws.Range("J74").Value = 22565.334
$ws.Range("L74").Value = 22565.334
$ws.Range("N74").Value = -24437.334

$ws.Range("H77").Value = 25003.5
$ws.Range("J77").Value = 22565.334
$ws.Range("L77").Value = 67696.00199999999
$ws.Range("N77").Value = -77056.00199999999

$ws.Range("H81").Value = 845070.5
$ws.Range("I81").Value = 8489.0625
$ws.Range("K81").Value = 16978.125
$ws.Range("M81").Value = -15917.125

$ws.Range("H84").Value = 845070.5
$ws.Range("I84").Value = 8489.0625
$ws.Range("K84").Value = 84890.625
$ws.Range("M84").Value = -79586.625

$ws.Range("H122").Value = 8120.6206
$ws.Range("J122").Value = 15187.637
$ws.Range("L122").Value = 45562.911
$ws.Range("N122").Value = -50462.911
